$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text block on sheet Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.67 = 6259.43 pesos`n✅ 6259.43 pesos = 1.66 = 897.94 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update tasas rates on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 597.5
$wsTasas.Range("O10").Value = 3740.01
$wsTasas.Range("N12").Value = 3764.99
$wsTasas.Range("O12").Value = 540.1
